# com soc project item list.xlsx - "Add files via upload" edit
# Updates the Raspberry Pi power-supply line item, turns the old
# "external power supply for servos" / "(servo team working on it)" row
# into a proper "servo controller" line with a real link + price, adds a
# new "servo controller power supply" row, and extends the total/footnote.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Raspberry Pi row: new variant link + new price
# ---------------------------------------------------------------------
$ws.Range("B2").Value = "https://thepihut.com/products/raspberry-pi-4-model-b?variant=20064052740158"
$ws.Range("C2").Value = 54

# ---------------------------------------------------------------------
# 2) Row 12 becomes "servo controller" (was the placeholder
#    "external power supply for servos" / "(servo team working on it)" row)
# ---------------------------------------------------------------------
$ws.Range("A12").Value = "servo controller"
$ws.Range("B12").Value = "https://www.adafruit.com/product/815 "
$ws.Range("C12").Value = 14.95

# ---------------------------------------------------------------------
# 3) New row 13: "servo controller power supply"
# ---------------------------------------------------------------------
$ws.Range("A13").Value = "servo controller power supply"
$ws.Range("B13").Value = "https://www.amazon.co.uk/Adapter-aifulo-Universal-Adapters-Speakers/dp/B08CZ9VTSH/ref=sr_1_3?dchild=1&keywords=2.1mm+DC+jack+5v+10A&qid=1601824635&sr=8-3 "
$ws.Range("C13").Value = 8.99

# ---------------------------------------------------------------------
# 4) Move the running total down to row 14 and add the budget footnote
# ---------------------------------------------------------------------
$ws.Range("C14").Formula = "=SUM(C2:C13)"
$ws.Range("D14").Value = "(max is £250)"

# ---------------------------------------------------------------------
# 5) Hyperlinks: this engine's Hyperlinks.Delete() always clears every
#    hyperlink on the sheet (range-scoped delete isn't respected), so the
#    only reliable way to update B2's target and add the two new ones is
#    to drop all of them and recreate the full set in place.
# ---------------------------------------------------------------------
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("B2"), "https://thepihut.com/products/raspberry-pi-4-model-b?variant=20064052740158")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://www.adafruit.com/product/2143")

$h4 = $ws.Hyperlinks.Add($ws.Range("B4"), "https://www.amazon.co.uk/AZDelivery-MB-102-Breadboard-Kit/dp/B07KYHBVR7/ref=sr_1_2_sspa?crid=3JE7OW6UECZDW&dchild=1&keywords=raspberry+pi+gpio+extension&qid=1593113124&sprefix=raspberry+pi+gpio+%2Caps%2C145&sr=8-2-spons&psc=1&spLa=ZW5jcnlwdGVkUXVhbGlmaWVyPUEzVERNQTE3STJUQ0hHJmVuY3J5cHRlZElkPUEwNDY1OTI0MVNRVzZGMzFGMjVLVCZlbmNyeXB0ZWRBZElkPUEwNDM2NDI5M1Q1TkkwUzFEMldOOSZ3aWRnZXROYW1lPXNwX2F0ZiZhY3Rpb249Y2xpY2tSZWRpcmVjdCZkb05vdExvZ0NsaWNrPXRydWU=")
$h4.TextToDisplay = "https://www.amazon.co.uk/AZDelivery-MB-102-Breadboard-Kit/dp/B07KYHBVR7/ref=sr_1_2_sspa?crid=3JE7OW6UECZDW&dchild=1&keywords=raspberry+pi+gpio+extension&qid=1593113124&sprefix=raspberry+pi+gpio+%2Caps%2C145&sr=8-2-spons&psc=1&spLa=ZW5jcnlwdGVkUXVhbGlmaWVyPUEzVERNQTE3STJUQ0hHJmVuY3J5cHRlZElkPUEwNDY1OTI0MVNRVzZGMzFGMjVLVCZlbmNyeXB0ZWRBZElkPUEwNDM2NDI5M1Q1TkkwUzFEMldOOSZ3aWRnZXROYW1lPXNwX2F0ZiZhY3Rpb249Y2xpY2tSZWRpcmVjdCZkb05vdExvZ0NsaWNrPXRydWU="

$h5 = $ws.Hyperlinks.Add($ws.Range("B5"), "https://www.amazon.co.uk/100-Pack-Black-Cable-Ties/dp/B07BGTR4G4/ref=sr_1_1_sspa?dchild=1&keywords=cable+ties&qid=1593113138&sr=8-1-spons&psc=1&spLa=ZW5jcnlwdGVkUXVhbGlmaWVyPUEzMzBITjZaQjg2TktSJmVuY3J5cHRlZElkPUEwMjQwMzc2MUtSRDdURDFDR0NUUyZlbmNyeXB0ZWRBZElkPUEwOTI1MTA0TjVNOU8wSE43SDlCJndpZGdldE5hbWU9c3BfYXRmJmFjdGlvbj1jbGlja1JlZGlyZWN0JmRvTm90TG9nQ2xpY2s9dHJ1ZQ==")
$h5.TextToDisplay = "https://www.amazon.co.uk/100-Pack-Black-Cable-Ties/dp/B07BGTR4G4/ref=sr_1_1_sspa?dchild=1&keywords=cable+ties&qid=1593113138&sr=8-1-spons&psc=1&spLa=ZW5jcnlwdGVkUXVhbGlmaWVyPUEzMzBITjZaQjg2TktSJmVuY3J5cHRlZElkPUEwMjQwMzc2MUtSRDdURDFDR0NUUyZlbmNyeXB0ZWRBZElkPUEwOTI1MTA0TjVNOU8wSE43SDlCJndpZGdldE5hbWU9c3BfYXRmJmFjdGlvbj1jbGlja1JlZGlyZWN0JmRvTm90TG9nQ2xpY2s9dHJ1ZQ=="

$ws.Hyperlinks.Add($ws.Range("B6"), "https://thepihut.com/products/official-raspberry-pi-universal-power-supply?src=raspberrypi")

$h7 = $ws.Hyperlinks.Add($ws.Range("B7"), "https://www.amazon.co.uk/ANYCUBIC-Printer-Filament-Printing-Printers/dp/B07DMF9ZRL/ref=alex_attr_sims_c_t1_3/259-3827962-2870757?_encoding=UTF8&pd_rd_i=B07DMF9ZRL&pd_rd_r=c703ee8c-572c-4524-a456-3c8c0d2822cb&pd_rd_w=Kifop&pd_rd_wg=jvnsN&pf_rd_p=5a049b35-22e7-4a5d-9586-cabe711af4ee&pf_rd_r=83ECAFVC33X25A4088MH&psc=1&refRID=83ECAFVC33X25A4088MH")
$h7.TextToDisplay = "https://www.amazon.co.uk/ANYCUBIC-Printer-Filament-Printing-Printers/dp/B07DMF9ZRL/ref=alex_attr_sims_c_t1_3/259-3827962-2870757?_encoding=UTF8&pd_rd_i=B07DMF9ZRL&pd_rd_r=c703ee8c-572c-4524-a456-3c8c0d2822cb&pd_rd_w=Kifop&pd_rd_wg=jvnsN&pf_rd_p=5a049b35-22e7-4a5d-9586-cabe711af4ee&pf_rd_r=83ECAFVC33X25A4088MH&psc=1&refRID=83ECAFVC33X25A4088MH"

$h8 = $ws.Hyperlinks.Add($ws.Range("B8"), "https://www.amazon.co.uk/Diymore-Digital-Helicopter-Airplane-controls/dp/B07DQJ1JXY/ref=sr_1_1_sspa?dchild=1&keywords=MG996R%2Bservo&qid=1593112182&sr=8-1-spons&spLa=ZW5jcnlwdGVkUXVhbGlmaWVyPUFJR1M3TEJPRTNXVkomZW5jcnlwdGVkSWQ9QTA0OTIyMjkxUVZPVU1URlJJUjJDJmVuY3J5cHRlZEFkSWQ9QTEwMzU3ODlYQ1hOTkU4OEU0N0cmd2lkZ2V0TmFtZT1zcF9hdGYmYWN0aW9uPWNsaWNrUmVkaXJlY3QmZG9Ob3RMb2dDbGljaz10cnVl&th=1")
$h8.TextToDisplay = "https://www.amazon.co.uk/Diymore-Digital-Helicopter-Airplane-controls/dp/B07DQJ1JXY/ref=sr_1_1_sspa?dchild=1&keywords=MG996R%2Bservo&qid=1593112182&sr=8-1-spons&spLa=ZW5jcnlwdGVkUXVhbGlmaWVyPUFJR1M3TEJPRTNXVkomZW5jcnlwdGVkSWQ9QTA0OTIyMjkxUVZPVU1URlJJUjJDJmVuY3J5cHRlZEFkSWQ9QTEwMzU3ODlYQ1hOTkU4OEU0N0cmd2lkZ2V0TmFtZT1zcF9hdGYmYWN0aW9uPWNsaWNrUmVkaXJlY3QmZG9Ob3RMb2dDbGljaz10cnVl&th=1"

$ws.Hyperlinks.Add($ws.Range("B9"), "https://thepihut.com/products/raspberry-pi-camera-module?src=raspberrypi")

$h10 = $ws.Hyperlinks.Add($ws.Range("B10"), "https://www.amazon.co.uk/Adhesive-Bumpers-Cylindrical-Hemispherical-Trapezoid/dp/B01M7OAWNU/ref=sr_1_1_sspa?dchild=1&keywords=rubber+sticky+pads&qid=1593116519&sr=8-1-spons&psc=1&spLa=ZW5jcnlwdGVkUXVhbGlmaWVyPUEzOUhTTExWOUpKSkg2JmVuY3J5cHRlZElkPUEwNDgxODkzMlNIQUtSUFlaN1dLOCZlbmNyeXB0ZWRBZElkPUEwODU2NDQwMkNQOUIxWFZUVjFaTyZ3aWRnZXROYW1lPXNwX2F0ZiZhY3Rpb249Y2xpY2tSZWRpcmVjdCZkb05vdExvZ0NsaWNrPXRydWU=")
$h10.TextToDisplay = "https://www.amazon.co.uk/Adhesive-Bumpers-Cylindrical-Hemispherical-Trapezoid/dp/B01M7OAWNU/ref=sr_1_1_sspa?dchild=1&keywords=rubber+sticky+pads&qid=1593116519&sr=8-1-spons&psc=1&spLa=ZW5jcnlwdGVkUXVhbGlmaWVyPUEzOUhTTExWOUpKSkg2JmVuY3J5cHRlZElkPUEwNDgxODkzMlNIQUtSUFlaN1dLOCZlbmNyeXB0ZWRBZElkPUEwODU2NDQwMkNQOUIxWFZUVjFaTyZ3aWRnZXROYW1lPXNwX2F0ZiZhY3Rpb249Y2xpY2tSZWRpcmVjdCZkb05vdExvZ0NsaWNrPXRydWU="

$ws.Hyperlinks.Add($ws.Range("B11"), "https://thepihut.com/products/usb-b-to-usb-c-adapter")

# New hyperlinks for the servo-controller rows
$ws.Hyperlinks.Add($ws.Range("B12"), "https://www.adafruit.com/product/815 ")
$ws.Hyperlinks.Add($ws.Range("B13"), "https://www.amazon.co.uk/Adapter-aifulo-Universal-Adapters-Speakers/dp/B08CZ9VTSH/ref=sr_1_3?dchild=1&keywords=2.1mm+DC+jack+5v+10A&qid=1601824635&sr=8-3 ")

# Re-apply the Hyperlink cell style everywhere a link lives (Hyperlinks.Add
# swaps in a throwaway style internally, so put the normal "Hyperlink"
# look back across the whole link column in one shot).
$ws.Range("B2:B13").Style = "Hyperlink"

# ---------------------------------------------------------------------
# 6) Selection, as last left by the editor
# ---------------------------------------------------------------------
$ws.Range("D17").Select()
